$wb = $excel.ActiveWorkbook

# --- Worksheets ---
$wsAbout = $wb.Worksheets.Item("About")
$wsCurrentPlanned = $wb.Worksheets.Item("Current and Planned Capacity")
$wsBAU = $wb.Worksheets.Item("BAU Emissions")

# --- Update the "NoSettings" suffix to "test" on every row label in BAU Emissions ---
$wsBAU.Cells.Replace(": NoSettings", ": test")

# --- Update the report date on the About sheet ---
$wsAbout.Range("C1").Value = 45387

# --- Update the BAU Emissions figures for row 94 (cols M:AE) ---
$cols = @("M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE")
$vals = @(1001080,2002150,3003230,4004300,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsBAU.Range($cols[$i] + "94").Value = $vals[$i]
}

# --- Update selection/scroll state on sheets whose view changed ---
$wsBAU.Range("A30:AE280").Select()

# --- Make "About" the active sheet/tab (must be done last so it becomes the selected tab) ---
$wsAbout.Activate()
$wsAbout.Range("E29").Select()
